$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date serial for each row.
# Update all data rows (2 through 357) from 45172 (2023-09-03) to 45175 (2023-09-06).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 357 }

$ws.Range("C2:C$lastRow").Value = 45175
